$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: updated timestamp string
$ws.Range("A1").Value = 'Datos actualizados a 1 de Mayo de 2020 a las 01:22'

# Data rows: country (col A) + 7 numeric stats (cols B-H)
# Caused by countries changing rank order in the daily update
$rows = @(
  @{ Row = 4; Country = 'Estados Unidos'; Vals = @(1094464, 30270, 151818, 878819, 15226, 2172, 63827) },
  @{ Row = 15; Country = 'Canada'; Vals = @(53236, 1639, 21423, 28629, 557, 188, 3184) },
  @{ Row = 45; Country = 'Chequia'; Vals = @(7682, 103, 3314, 4132, 68, 9, 236) },
  @{ Row = 46; Country = 'Banglades'; Vals = @(7667, 564, 160, 7339, 1, 5, 168) },
  @{ Row = 49; Country = 'Colombia'; Vals = @(6507, 300, 1439, 4775, 118, 15, 293) },
  @{ Row = 50; Country = 'Panama'; Vals = @(6378, 178, 527, 5673, 92, 2, 178) },
  @{ Row = 73; Country = 'Nigeria'; Vals = @(1932, 204, 319, 1555, 2, 7, 58) },
  @{ Row = 74; Country = 'Camerun'; Vals = @(1832, 0, 934, 837, 12, 0, 61) },
  @{ Row = 75; Country = 'Azerbaiyan'; Vals = @(1804, 38, 1325, 455, 17, 1, 24) },
  @{ Row = 76; Country = 'Islandia'; Vals = @(1797, 0, 1670, 117, 0, 0, 10) },
  @{ Row = 77; Country = 'Bosnia y Herzegovina'; Vals = @(1757, 80, 727, 961, 4, 4, 69) },
  @{ Row = 105; Country = 'Uruguay'; Vals = @(643, 18, 417, 209, 10, 2, 17) },
  @{ Row = 156; Country = 'Guyana'; Vals = @(82, 8, 22, 51, 2, 1, 9) },
  @{ Row = 157; Country = 'Liechtenstein'; Vals = @(82, 0, 55, 26, 0, 0, 1) },
  @{ Row = 158; Country = 'Barbados'; Vals = @(81, 1, 39, 35, 4, 0, 7) },
  @{ Row = 159; Country = 'Bahamas'; Vals = @(80, 0, 25, 44, 1, 0, 8) }
)

foreach ($r in $rows) {
  $ws.Cells.Item($r.Row, 1).Value = $r.Country
  $col = 2
  foreach ($v in $r.Vals) {
    $ws.Cells.Item($r.Row, $col).Value = $v
    $col = $col + 1
  }
}
